$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 19-42: station name (col A), terminal name (col B), last-charge-end timestamp (col C)
# get reshuffled/updated in this edit (new scrape results appended, old ones trimmed).
$rows = @(
    @{ Row = 19; A = "长沙特来电飞狐四方坪南区充电站";             B = "406号直流";  C = 45943.020914351851 },
    @{ Row = 20; A = "长沙特来电飞狐四方坪南区充电站";             B = "306号直流";  C = 45944.674907407411 },
    @{ Row = 21; A = "长沙特来电飞狐四方坪东区充电站";             B = "904号直流";  C = 45945.536215277774 },
    @{ Row = 22; A = "长沙市开福区高岭香江国际城充电站建设项目";   B = "109号直流";  C = 45945.603726851848 },
    @{ Row = 23; A = "长沙特来电飞狐四方坪西区充电站";             B = "603号直流";  C = 45946.045289351852 },
    @{ Row = 24; A = "长沙特来电飞狐四方坪西区充电站";             B = "604号直流";  C = 45946.1093287037 },
    @{ Row = 25; A = "长沙特来电飞狐四方坪东区充电站";             B = "004A号直流"; C = 45946.16909722222 },
    @{ Row = 26; A = "长沙特来电飞狐四方坪西区充电站";             B = "703号直流";  C = 45946.185879629629 },
    @{ Row = 27; A = "长沙特来电飞狐四方坪东区充电站";             B = "011A号直流"; C = 45946.228773148148 },
    @{ Row = 28; A = "长沙特来电飞狐四方坪西区充电站";             B = "802号直流";  C = 45946.245138888888 },
    @{ Row = 29; A = "长沙特来电飞狐四方坪东区充电站";             B = "401号直流";  C = 45946.255624999998 },
    @{ Row = 30; A = "长沙特来电飞狐四方坪西区充电站";             B = "A01号直流";  C = 45946.297500000001 },
    @{ Row = 31; A = "长沙特来电飞狐四方坪西区充电站";             B = "401号直流";  C = 45946.298530092594 },
    @{ Row = 32; A = "长沙特来电飞狐四方坪西区充电站";             B = "904号直流";  C = 45946.326631944445 },
    @{ Row = 33; A = "长沙市开福区高岭香江国际城充电站建设项目";   B = "107号直流";  C = 45946.35597222222 },
    @{ Row = 34; A = "长沙市开福区高岭香江国际城充电站建设项目";   B = "207号直流";  C = 45946.411979166667 },
    @{ Row = 35; A = "长沙特来电飞狐四方坪南区充电站";             B = "305号直流";  C = 45946.441458333335 },
    @{ Row = 36; A = "长沙市开福区高岭香江国际城充电站建设项目";   B = "108号直流";  C = 45946.510474537034 },
    @{ Row = 37; A = "长沙特来电飞狐四方坪南区充电站";             B = "206号直流";  C = 45946.545810185184 },
    @{ Row = 38; A = "长沙特来电飞狐四方坪西区充电站";             B = "902号直流";  C = 45946.547372685185 },
    @{ Row = 39; A = "长沙特来电飞狐四方坪西区充电站";             B = "602号直流";  C = 45946.549641203703 },
    @{ Row = 40; A = "长沙特来电飞狐四方坪南区充电站";             B = "401号直流";  C = 45946.551319444443 },
    @{ Row = 41; A = "长沙特来电飞狐四方坪南区充电站";             B = "106号直流";  C = 45946.554664351854 },
    @{ Row = 42; A = "长沙特来电飞狐四方坪西区充电站";             B = "505号直流";  C = 45946.557071759256 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value2 = $r.A
    $ws.Cells.Item($r.Row, 2).Value2 = $r.B
    $ws.Cells.Item($r.Row, 3).Value2 = $r.C
}

# Rows 43-48 no longer have data in this refresh - clear A:C (D/E were already blank).
for ($row = 43; $row -le 48; $row++) {
    $ws.Range("A$row`:C$row").ClearContents()
}

# Restore the saved selection/active cell from the workbook.
$ws.Range("E23").Select()
